$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_R_acc_LT"

$values = @(
    88.586479367866545,
    91.483757682177341,
    91.17647058823529,
    81.60667251975417,
    87.093942054433711,
    85.425812115891134,
    91.439859525899919,
    91.220368744512726,
    91.527655838454777,
    91.791044776119406,
    76.119402985074629,
    77.216856892010526,
    88.542581211589109,
    88.84986830553116,
    88.498683055311673,
    71.773485513608421,
    73.222124670763833,
    73.617208077260756,
    88.805970149253739,
    89.069359086918354,
    89.947322212467071,
    76.031606672519757,
    75.373134328358205,
    75.592625109745399,
    74.670763827919231,
    73.266022827041269,
    74.539069359086923,
    88.937664618086046,
    88.366988586479366,
    73.924495171202807,
    87.006145741878839,
    91.878841088674278,
    89.157155399473226,
    89.201053555750647,
    90.342405618964008,
    67.515364354697098,
    87.928007023705007,
    86.128182616330122,
    85.645302897278313,
    83.230904302019312,
    72.519754170324845,
    73.485513608428448,
    83.494293239683941,
    84.635645302897274,
    83.187006145741876,
    79.587357330992106,
    73.090430201931525,
    81.738366988586478,
    79.806848112379285
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
